# Updated cryptos list values (Price / Volume(1h)) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.302.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.872.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7099'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.74'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07806'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08400'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.34'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.09'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.312.42'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.072'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008180'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.97'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.22'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.121.61'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.762'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.15'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.996'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.45'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.505'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.391'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.296'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.295'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05387'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.943'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7501'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.693'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01869'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.231.12'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.728'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.529'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8885'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.44'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.53'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.021.96'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5194'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.422'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.53%  '
